# Allow importing new responses without an initial response (#1543)
#
# Rename the first sheet and swap which sheet/cell is active & selected:
#   - Sheet1 "Basic Clinic Data..." -> "Test Survey", becomes the active
#     (tabSelected) sheet with E1 selected.
#   - Sheet2 "Facility Fundamen..." is no longer the active sheet; its
#     stored selection becomes G7.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename the first worksheet.
$ws1.Name = "Test Survey"

# Update sheet2's saved selection first (while it's still active), then
# hand activation over to sheet1 so it ends up as the tabSelected sheet.
$ws2.Activate()
$ws2.Range("G7").Select() | Out-Null

$ws1.Activate()
$ws1.Range("E1").Select() | Out-Null
